# SectorGroup.xlsx — reorder the codeforiati group/category columns.
#
# The sheet has columns:
#   A=code, B=name, C=status,
#   D=codeforiati:group-code, E=codeforiati:category-name,
#   F=codeforiati:category-code, G=codeforiati:group-name
#
# The edit relabels/reorders the last four columns to:
#   D=codeforiati:category-name, E=codeforiati:group-name,
#   F=codeforiati:category-code, G=codeforiati:group-code
#
# i.e. a 3-way rotation of the D/E/G values on every data row
# (new D = old E, new E = old G, new G = old D); F (category-code)
# keeps its value. The header row is relabeled the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# --- Header row -------------------------------------------------------
$ws.Range("D1").Value = "codeforiati:category-name"
$ws.Range("E1").Value = "codeforiati:group-name"
$ws.Range("F1").Value = "codeforiati:category-code"
$ws.Range("G1").Value = "codeforiati:group-code"

# --- Data rows ----------------------------------------------------------
# Force text storage on D/E/G so numeric-looking codes ("110", "230", ...)
# stay strings instead of being coerced to numbers.
$ws.Range("D2:D$lastRow").NumberFormat = "@"
$ws.Range("E2:E$lastRow").NumberFormat = "@"
$ws.Range("G2:G$lastRow").NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $oldD = $ws.Cells.Item($r, 4).Value2
    $oldE = $ws.Cells.Item($r, 5).Value2
    $oldG = $ws.Cells.Item($r, 7).Value2

    $ws.Cells.Item($r, 4).Value = $oldE
    $ws.Cells.Item($r, 5).Value = $oldG
    $ws.Cells.Item($r, 7).Value = $oldD
}
